# Update the two-digit division worksheet numbers.
# Each (old, new) pair below corresponds to one "<dividend>÷<divisor>=" cell
# that changed between the original and the committed version, in the
# same order they occur in the document.

$d = $word.ActiveDocument

$pairs = @(
    @("36÷9=", "30÷6="),
    @("30÷2=", "23÷7="),
    @("72÷3=", "92÷2="),
    @("23÷9=", "26÷4="),
    @("24÷2=", "67÷4="),
    @("88÷5=", "52÷2="),
    @("33÷3=", "86÷5="),
    @("70÷9=", "53÷7="),
    @("84÷5=", "58÷8="),
    @("32÷9=", "61÷8="),
    @("66÷7=", "70÷7="),
    @("39÷4=", "37÷9="),
    @("87÷6=", "35÷2="),
    @("90÷5=", "83÷8="),
    @("53÷8=", "11÷7="),
    @("85÷4=", "62÷2="),
    @("93÷7=", "68÷8="),
    @("56÷8=", "52÷2="),
    @("71÷4=", "91÷3="),
    @("35÷7=", "88÷7="),
    @("44÷5=", "47÷5="),
    @("80÷7=", "84÷3="),
    @("29÷3=", "58÷5="),
    @("94÷9=", "68÷7="),
    @("79÷3=", "35÷7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true,
                         1, $false, $new, 2)
}
